# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  "Office Theme" (clrScheme "Office")      <->
#   ppt/theme/theme2.xml  "Integral"     (clrScheme "Red Violet")
#
# ppt/theme/theme2.xml is the theme actually bound to the slide master
# (and to the presentation itself), so it is the part that drives what
# colors the deck's slides render with. After the swap it must carry the
# "Office Theme" color values that theme1.xml used to hold (name/font
# scheme/format scheme are identical between the two themes already, so
# the 12 scheme colors are the only observable difference).
#
# RGB(r,g,b) == r + g*256 + b*65536, matching VBA's color-long encoding.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink 954F72
